$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '36.926.19'
$ws.Range("E2").Value = '  +2.31%  '
$ws.Range("D3").Value = '2.036.26'
$ws.Range("E3").Value = '  +1.43%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.96'
$ws.Range("E5").Value = '  -1.06%  '
$ws.Range("E6").Value = '  -0.62%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '62.95'
$ws.Range("E7").Value = '  +1.49%  '
$ws.Range("E8").Value = '  +0.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.388'
$ws.Range("E9").Value = '  +5.53%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '58.20'
$ws.Range("E10").Value = '  -1.32%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0797'
$ws.Range("E11").Value = '  +7.01%  '
$ws.Range("E12").Value = '  -0.68%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.906'
$ws.Range("E13").Value = '  -1.31%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '23.12'
$ws.Range("E14").Value = '  +17.99%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '14.37'
$ws.Range("E15").Value = '  -3.03%  '
$ws.Range("D16").Value = '2.336.01'
$ws.Range("E16").Value = '  +1.44%  '
$ws.Range("E17").Value = '  +2.65%  '
$ws.Range("D18").Value = '2.036.39'
$ws.Range("E18").Value = '  +1.53%  '
$ws.Range("D19").Value = '36.896.94'
$ws.Range("E19").Value = '  +2.39%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '72.24'
$ws.Range("E20").Value = '  +0.48%  '
$ws.Range("D21").Value = '0.0₃0883'
$ws.Range("E21").Value = '  +3.53%  '
$ws.Range("E22").Value = '  +2.80%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '236.36'
$ws.Range("E23").Value = '  +1.32%  '
$ws.Range("E24").Value = '  +0.17%  '
$ws.Range("E25").Value = '  -7.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.33'
$ws.Range("E26").Value = '  +0.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.77'
$ws.Range("E27").Value = '  +3.12%  '
$ws.Range("B28").Value = 'Kaspa'
$ws.Range("C28").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.143'
$ws.Range("E28").Value = '  +33.05%  '
$ws.Range("B29").Value = 'Monero'
$ws.Range("C29").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '159.57'
$ws.Range("E29").Value = '  -3.01%  '
$ws.Range("B30").Value = 'EthereumClassic'
$ws.Range("C30").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '20.23'
$ws.Range("E30").Value = '  +3.58%  '
$ws.Range("E31").Value = '  +1.23%  '
$ws.Range("E32").Value = '  -0.95%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.18'
$ws.Range("E33").Value = '  -0.73%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0620'
$ws.Range("E34").Value = '  +3.03%  '
$ws.Range("E35").Value = '  +1.01%  '
$ws.Range("E36").Value = '  -3.54%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '6.32'
$ws.Range("E37").Value = '  +10.35%  '
$ws.Range("E38").Value = '  -0.09%  '
$ws.Range("E39").Value = '  +1.88%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.13'
$ws.Range("E40").Value = '  +33.68%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.0991'
$ws.Range("E41").Value = '  -9.45%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.25'
$ws.Range("E42").Value = '  +2.70%  '
$ws.Range("E43").Value = '  +3.32%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '17.10'
$ws.Range("E44").Value = '  +2.07%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.14'
$ws.Range("E45").Value = '  +1.71%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0215'
$ws.Range("E46").Value = '  +0.04%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '93.64'
$ws.Range("E47").Value = '  +0.37%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.73'
$ws.Range("E48").Value = '  -1.01%  '
$ws.Range("D49").Value = '1.368.99'
$ws.Range("E49").Value = '  -3.01%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.92'
$ws.Range("E50").Value = '  +0.62%  '
$ws.Range("D51").Value = '2.222.72'
$ws.Range("E51").Value = '  +1.36%  '
